# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" row at the top of the "总计" (summary) sheet,
#     pushing the existing quarterly rows down by one.
#  2. Insert a brand-new "2022-Q4" worksheet (cloned from the existing
#     "2022-Q3" sheet so it keeps the same header/styling) right before
#     the "2022-Q3" tab, and fill it with the per-fund holdings data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift quarters down one row and insert 2022-Q4 on top
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift columns B (date label), C (count) and D (market value) down by
# one row, from the bottom up so we don't clobber values before reading
# them. Column A is just a 0-based running index, so it's rebuilt below
# rather than shifted.
for ($r = 8; $r -ge 2; $r--) {
    $total.Cells.Item($r + 1, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($r + 1, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($r + 1, 4).Value = $total.Cells.Item($r, 4).Value2
}

# New top row: 2022-Q4 totals
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 10
$total.Cells.Item(2, 4).Value = 0.48

# Extend column A's styled run down to the new row 9 (copies style + value
# from row 8, value gets overwritten right after).
$total.Cells.Item(8, 1).Copy($total.Cells.Item(9, 1))

# Rebuild column A as the sequential index 0..7 for rows 2..9
for ($r = 2; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with fund holdings, placed right before the
#    existing "2022-Q3" tab.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
# NOTE: after Copy(Before) inserts the clone in front of $q3, this engine's
# `$q3` handle keeps tracking the *position* the variable was bound to
# rather than the original sheet object - so it now refers to the newly
# inserted clone (which sits at that same slot), while the untouched
# original "2022-Q3" sheet has shifted one slot later. That's exactly the
# sheet we want to turn into "2022-Q4", so just keep using $q3.
$ws = $q3
$ws.Name = "2022-Q4"

# The source sheet (2022-Q3) only had 6 data rows; 2022-Q4 needs 10, so
# stamp out 4 more rows with the same formatting as row 7 before writing
# the real values over everything.
for ($r = 8; $r -le 11; $r++) {
    $ws.Range("A7:H7").Copy($ws.Range("A" + $r + ":H" + $r))
}

# row, code, name, scale, stock position, position ratio, holding value, rank
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'378010"
$ws.Cells.Item(2,3).Value = "上投摩根成长先锋混合A"
$ws.Cells.Item(2,4).Value = "'9.92"
$ws.Cells.Item(2,5).Value = "'88.07"
$ws.Cells.Item(2,6).Value = "'2.64"
$ws.Cells.Item(2,7).Value = "'0.2619"
$ws.Cells.Item(2,8).Value = 8

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'370024"
$ws.Cells.Item(3,3).Value = "上投摩根核心优选混合A"
$ws.Cells.Item(3,4).Value = "'7.47"
$ws.Cells.Item(3,5).Value = "'77.24"
$ws.Cells.Item(3,6).Value = "'2.13"
$ws.Cells.Item(3,7).Value = "'0.1591"
$ws.Cells.Item(3,8).Value = 10

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'004194"
$ws.Cells.Item(4,3).Value = "招商中证1000指数增强A"
$ws.Cells.Item(4,4).Value = "'2.57"
$ws.Cells.Item(4,5).Value = "'94.27"
$ws.Cells.Item(4,6).Value = "'1.12"
$ws.Cells.Item(4,7).Value = "'0.0288"
$ws.Cells.Item(4,8).Value = 5

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'004195"
$ws.Cells.Item(5,3).Value = "招商中证1000指数增强C"
$ws.Cells.Item(5,4).Value = "'2.14"
$ws.Cells.Item(5,5).Value = "'94.27"
$ws.Cells.Item(5,6).Value = "'1.12"
$ws.Cells.Item(5,7).Value = "'0.0240"
$ws.Cells.Item(5,8).Value = 5

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'011054"
$ws.Cells.Item(6,3).Value = "申万菱信安鑫智选混合A"
$ws.Cells.Item(6,4).Value = "'0.69"
$ws.Cells.Item(6,5).Value = "'25.20"
$ws.Cells.Item(6,6).Value = "'1.07"
$ws.Cells.Item(6,7).Value = "'0.0074"
$ws.Cells.Item(6,8).Value = 8

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'519222"
$ws.Cells.Item(7,3).Value = "海富通欣益灵活配置混合A"
$ws.Cells.Item(7,4).Value = "'0.25"
$ws.Cells.Item(7,5).Value = "'31.65"
$ws.Cells.Item(7,6).Value = "'0.16"
$ws.Cells.Item(7,7).Value = "'0.0004"
$ws.Cells.Item(7,8).Value = 10

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'015057"
$ws.Cells.Item(8,3).Value = "上投摩根核心优选混合C"
$ws.Cells.Item(8,4).Value = "'0.01"
$ws.Cells.Item(8,5).Value = "'77.24"
$ws.Cells.Item(8,6).Value = "'2.13"
$ws.Cells.Item(8,7).Value = "'0.0002"
$ws.Cells.Item(8,8).Value = 10

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'519221"
$ws.Cells.Item(9,3).Value = "海富通欣益灵活配置混合C"
$ws.Cells.Item(9,4).Value = "'0.10"
$ws.Cells.Item(9,5).Value = "'31.65"
$ws.Cells.Item(9,6).Value = "'0.16"
$ws.Cells.Item(9,7).Value = "'0.0002"
$ws.Cells.Item(9,8).Value = 10

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'015077"
$ws.Cells.Item(10,3).Value = "上投摩根成长先锋混合C"
$ws.Cells.Item(10,4).Value = "'0.00"
$ws.Cells.Item(10,5).Value = "'88.07"
$ws.Cells.Item(10,6).Value = "'2.64"
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 8

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'011055"
$ws.Cells.Item(11,3).Value = "申万菱信安鑫智选混合C"
$ws.Cells.Item(11,4).Value = "'0.00"
$ws.Cells.Item(11,5).Value = "'25.20"
$ws.Cells.Item(11,6).Value = "'1.07"
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 8

# The leading apostrophes above force text storage (the fund codes/ratios
# look numeric but must round-trip as plain strings), but they also tag
# the cells with a "quote prefix" style. Strip that back off now that the
# values are safely stored as text, restoring the plain/default style the
# source data used (columns B, D, E and F carry forced-text values; G is
# mixed - numeric 0 in rows 10-11 - so only clean its text rows).
$ws.Range("B2:B11").Style = "Normal"
$ws.Range("D2:F11").Style = "Normal"
$ws.Range("G2:G9").Style = "Normal"
